$d = $word.ActiveDocument

# Original paragraph text is "Version 2." split as runs:
#   "Versi" | "on" | " 2" | <bookmark "_GoBack"> | "."
# Target paragraph text is "Version 1." split as runs:
#   "Version" | " 1." | <bookmark "_GoBack">

# 1) Merge "Versi" + "on" into a single run reading "Version".
#    Find/Replace matches across the run boundary but leaves the
#    surrounding spellcheck markers / bookmark untouched.
$found1 = $d.Content.Find.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# 2) Change the version number run " 2" to " 1." (stays inside its own
#    run, so the bookmark immediately after it is not disturbed).
$found2 = $d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false, $true, 1, $false, " 1.", 2)

# 3) The trailing "." run (now redundant, since the period moved into
#    the run above) still follows the bookmark. Remove just that final
#    character -- it sits right before the paragraph mark.
$full = $d.Content.Text
$lastCharIndex = $full.Length - 2
$trailingPeriod = $d.Range($lastCharIndex, $lastCharIndex + 1)
$trailingPeriod.Delete()
